{"js": "// Append \" (Changed main)\" to the end of the sentence \"This is a\n// Microsoft word document.\", as three additional runs (\" (\",\n// \"Changed main\", \")\") sitting alongside the original run \u2014 rather\n// than growing the original run's text.\n//\n// A plain insertText() call would get silently coalesced into the\n// existing run when the package is serialized (adjacent runs with\n// identical formatting merge), so instead we splice in literal OOXML\n// for the new runs via insertOoxml(), which preserves them as\n// independent <w:r> siblings \u2014 matching the target diff exactly.\n\nconst results = context.document.body.search(\n  \"This is a Microsoft word document.\",\n  { matchCase: true }\n);\nresults.load(\"items\");\nawait context.sync();\n\nconst sentence = results.items[0];\nconst endOfSentence = sentence.getRange(Word.RangeLocation.end);\n\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>' +\n  '<w:r><w:t>Changed main</w:t></w:r>' +\n  '<w:r><w:t>)</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nendOfSentence.insertOoxml(ooxml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the exact sentence via Find so this does not depend on the\n# paragraph's ordinal position in the document.\n$searchRange = $d.Content\n[void]$searchRange.Find.Execute(\"This is a Microsoft word document.\")\n\n$targetParagraph = $searchRange.Paragraphs(1)\n$paragraphRange = $targetParagraph.Range\n\n# Plain Range.InsertAfter() text gets silently coalesced into the\n# existing run on save when formatting matches, so instead we rebuild\n# the whole paragraph as literal OOXML (keeping its own w14:paraId /\n# w14:textId / rsid attributes via Range.WordOpenXML) with the three new\n# runs (\" (\", \"Changed main\", \")\") appended as independent <w:r>\n# siblings after the original run \u2014 matching the target diff exactly.\n$fullPackageXml = $paragraphRange.WordOpenXML\n\n$openParaTag = '<w:p>'\nif ($fullPackageXml -match '(?s)<w:body>(<w:p\\b[^>]*>)') {\n    $openParaTag = $matches[1]\n}\n\n$newParagraphXml = $openParaTag +\n    '<w:r><w:t>This is a Microsoft word document.</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>' +\n    '<w:r><w:t>Changed main</w:t></w:r>' +\n    '<w:r><w:t>)</w:t></w:r>' +\n    '</w:p>'\n\n$wrappedXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n    '<w:body>' + $newParagraphXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n\n[void]$paragraphRange.InsertXML($wrappedXml)\n"}
